# Update crypto price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (e.g. "66.845.66" dot-grouped, or
# "577.21"). Excel's COM layer auto-converts a numeric-looking string typed
# into a General-formatted cell into a real number, which would corrupt
# values such as "577.04" with binary floating-point noise. Force those
# specific cells to Text format first so the new price strings round-trip
# exactly; cells whose new value is not a bare number (it still has a dot
# thousands-separator) do not need this and are left as-is.
$riskyCells = @("D5","D6","D10","D14","D18","D20","D21","D23","D24","D25","D27","D29","D30","D32","D37","D38","D39","D40","D45","D48","D50")
foreach ($addr in $riskyCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '66.823.78'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '3.111.94'
$ws.Range("E3").Value = '  +0.69%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '577.04'
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("D6").Value = '171.01'
$ws.Range("E6").Value = '  +2.13%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '3.108.01'
$ws.Range("E8").Value = '  +0.69%  '

$ws.Range("E9").Value = '  -0.68%  '

$ws.Range("D10").Value = '6.46'

$ws.Range("E11").Value = '  -1.22%  '

$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("E13").Value = '  -1.77%  '

$ws.Range("D14").Value = '37.16'
$ws.Range("E14").Value = '  +1.26%  '

$ws.Range("E15").Value = '  -1.08%  '

$ws.Range("D16").Value = '3.630.32'
$ws.Range("E16").Value = '  +0.68%  '

$ws.Range("D17").Value = '66.808.91'
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").Value = '7.15'
$ws.Range("E18").Value = '  -0.95%  '

$ws.Range("D19").Value = '3.112.52'
$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("D20").Value = '16.31'
$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("D21").Value = '475.99'
$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").Value = '7.94'
$ws.Range("E23").Value = '  +5.10%  '

$ws.Range("D24").Value = '13.37'
$ws.Range("E24").Value = '  +4.16%  '

$ws.Range("D25").Value = '84.06'
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("E26").Value = '  -2.89%  '

$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").Value = '7.86'
$ws.Range("E29").Value = '  -2.39%  '

$ws.Range("D30").Value = '2.36'
$ws.Range("E30").Value = '  -2.25%  '

$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("D32").Value = '28.51'
$ws.Range("E32").Value = '  +1.12%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("D34").Value = '0.0₃0940'
$ws.Range("E34").Value = '  -8.68%  '

$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("D37").Value = '0.974'
$ws.Range("E37").Value = '  -2.98%  '

$ws.Range("D38").Value = '46.92'
$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("D39").Value = '2.06'
$ws.Range("E39").Value = '  -2.11%  '

$ws.Range("D40").Value = '50.04'
$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("E41").Value = '  -1.83%  '

$ws.Range("E42").Value = '  -0.65%  '

$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").Value = '2.839.48'
$ws.Range("E44").Value = '  +2.30%  '

$ws.Range("D45").Value = '384.56'
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("E46").Value = '  -1.51%  '

$ws.Range("E47").Value = '  -9.47%  '

$ws.Range("D48").Value = '135.78'
$ws.Range("E48").Value = '  +0.68%  '

$ws.Range("E49").Value = '  +0.00%  '

$ws.Range("D50").Value = '24.85'
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("E51").Value = '  -2.06%  '
